$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.5185682773590088
$ws.Range("B1").Value = 0.6510264277458191
$ws.Range("C1").Value = 0.9419897198677063
$ws.Range("D1").Value = 3.975461006164551
$ws.Range("E1").Value = 4.130427837371826
